$d = $word.ActiveDocument

# Locate the bold "Outros adendos:" heading paragraph. The commit adds a new,
# empty paragraph (carrying bold paragraph-mark formatting) immediately
# before it, right after the paragraph that ends with
# "...calculadora financeira.".
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Outros adendos*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the 'Outros adendos:' paragraph to anchor the new paragraph before."
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$insertPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)

# Insert a brand-new, completely empty paragraph (no runs) whose paragraph
# mark carries bold/bold-complex-script formatting, via a raw OOXML
# fragment so no placeholder run is synthesized.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

[void]$insertPoint.InsertXML($xml)

Write-Host "Inserted new empty bold paragraph before paragraph index" $targetIndex
